# Family feature, longitudnal analysis
# Append a new data row (row 24) to the test-reports sheet, extending the
# used range from A1:AW23 to A1:AW24.
#
# Notes on technique:
#  - Plain text values (and numbers) can be assigned directly via .Value.
#  - "2025-09-08" and "1" look like a date / number to the COM layer's
#    auto-detection, so they are forced to text (NumberFormat "@") before
#    assignment, then the temporary number format is stripped again with
#    ClearFormats() so the cell keeps its default (unstyled) appearance
#    while remaining text-typed.
#  - Several columns are blank in the source row but still need to exist
#    as (empty) text cells, matching the rest of the sheet's blank cells.
#    Assigning "" is a no-op in this COM layer, so a lone "'" (empty
#    quote-prefixed text) is written instead and then the formatting is
#    cleared the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24

function Set-TextValue($col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

function Set-EmptyTextValue($col) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'"
    $c.ClearFormats()
}

# A24 - Date (looks like a date, must stay literal text)
Set-TextValue 1 "2025-09-08"

# B24 - Report Type
$ws.Cells.Item($row, 2).Value = "Liver Function Test (LFT)"

# C24:F24 - Hemoglobin, RBC, WBC, Platelets
$ws.Cells.Item($row, 3).Value = 11.3
$ws.Cells.Item($row, 4).Value = 3.78
$ws.Cells.Item($row, 5).Value = 10.84
$ws.Cells.Item($row, 6).Value = 486

# G24, H24 - Glucose, Cholesterol (blank)
Set-EmptyTextValue 7
Set-EmptyTextValue 8

# I24, J24 - Blood Pressure Systolic/Diastolic
$ws.Cells.Item($row, 9).Value = 13
$ws.Cells.Item($row, 10).Value = 9

# K24:M24 - Heart Rate, Temperature, Notes (blank)
Set-EmptyTextValue 11
Set-EmptyTextValue 12
Set-EmptyTextValue 13

# N24:P24 - Total/Conjugated/Unconjugated Bilirubin
$ws.Cells.Item($row, 14).Value = 0.4
$ws.Cells.Item($row, 15).Value = 0.2
$ws.Cells.Item($row, 16).Value = 0.2

# Q24:W24 - SGOT, SGPT, Alkaline Phosphatase, Total Protein, Albumin, Globulin, A/G Ratio
$ws.Cells.Item($row, 17).Value = 34
$ws.Cells.Item($row, 18).Value = 27
$ws.Cells.Item($row, 19).Value = 360
$ws.Cells.Item($row, 20).Value = 6.5
$ws.Cells.Item($row, 21).Value = 46
$ws.Cells.Item($row, 22).Value = 1.9
$ws.Cells.Item($row, 23).Value = 2.4

# X24 - PCV/HCT (blank)
Set-EmptyTextValue 24

# Y24:AC24 - MCV, MCH, MCHC, RDW-CV, MPV
$ws.Cells.Item($row, 25).Value = 78
$ws.Cells.Item($row, 26).Value = 25.4
$ws.Cells.Item($row, 27).Value = 32.6
$ws.Cells.Item($row, 28).Value = 12.7
$ws.Cells.Item($row, 29).Value = 8.199999999999999

# AD24:AF24 - T3, T4, TSH
$ws.Cells.Item($row, 30).Value = 3
$ws.Cells.Item($row, 31).Value = 4
$ws.Cells.Item($row, 32).Value = 2.98

# AG24:AJ24 - Neutrophils, Lymphocytes, Monocytes, Eosinophils
$ws.Cells.Item($row, 33).Value = 26
$ws.Cells.Item($row, 34).Value = 70
$ws.Cells.Item($row, 35).Value = 2
$ws.Cells.Item($row, 36).Value = 2

# AK24 - Gamma Glutamyl Transferase (blank)
Set-EmptyTextValue 37

# AL24 - Patient Name
$ws.Cells.Item($row, 38).Value = "KASHV"

# AM24 - Patient Age (stored as text "1")
Set-TextValue 39 "1"

# AN24:AW24 - Patient Gender, Liver Size, Gall Bladder Status, Spleen Size,
# Pancreas Status, Right/Left Kidney Size, Urinary Bladder Status,
# Ultrasound Findings, Ultrasound Impression (all blank)
Set-EmptyTextValue 40
Set-EmptyTextValue 41
Set-EmptyTextValue 42
Set-EmptyTextValue 43
Set-EmptyTextValue 44
Set-EmptyTextValue 45
Set-EmptyTextValue 46
Set-EmptyTextValue 47
Set-EmptyTextValue 48
Set-EmptyTextValue 49
